$d = $word.ActiveDocument

$d.Content.Find.Execute("609×2=", $true, $false, $false, $false, $false, $true, 1, $false, "399×2=", 2)
$d.Content.Find.Execute("850×6=", $true, $false, $false, $false, $false, $true, 1, $false, "888×3=", 2)
$d.Content.Find.Execute("117×9=", $true, $false, $false, $false, $false, $true, 1, $false, "667×6=", 2)
$d.Content.Find.Execute("411×7=", $true, $false, $false, $false, $false, $true, 1, $false, "957×5=", 2)
$d.Content.Find.Execute("490×5=", $true, $false, $false, $false, $false, $true, 1, $false, "775×6=", 2)
$d.Content.Find.Execute("546×9=", $true, $false, $false, $false, $false, $true, 1, $false, "555×2=", 2)
$d.Content.Find.Execute("460×6=", $true, $false, $false, $false, $false, $true, 1, $false, "444×9=", 2)
$d.Content.Find.Execute("868×5=", $true, $false, $false, $false, $false, $true, 1, $false, "162×4=", 2)
$d.Content.Find.Execute("638×5=", $true, $false, $false, $false, $false, $true, 1, $false, "972×2=", 2)
$d.Content.Find.Execute("520×3=", $true, $false, $false, $false, $false, $true, 1, $false, "469×9=", 2)
$d.Content.Find.Execute("417×6=", $true, $false, $false, $false, $false, $true, 1, $false, "312×9=", 2)
$d.Content.Find.Execute("625×2=", $true, $false, $false, $false, $false, $true, 1, $false, "870×4=", 2)
$d.Content.Find.Execute("683×9=", $true, $false, $false, $false, $false, $true, 1, $false, "816×2=", 2)
$d.Content.Find.Execute("612×5=", $true, $false, $false, $false, $false, $true, 1, $false, "646×5=", 2)
$d.Content.Find.Execute("883×8=", $true, $false, $false, $false, $false, $true, 1, $false, "156×7=", 2)
$d.Content.Find.Execute("248×5=", $true, $false, $false, $false, $false, $true, 1, $false, "264×9=", 2)
$d.Content.Find.Execute("403×8=", $true, $false, $false, $false, $false, $true, 1, $false, "786×3=", 2)
$d.Content.Find.Execute("776×7=", $true, $false, $false, $false, $false, $true, 1, $false, "592×8=", 2)
$d.Content.Find.Execute("728×9=", $true, $false, $false, $false, $false, $true, 1, $false, "781×7=", 2)
$d.Content.Find.Execute("642×2=", $true, $false, $false, $false, $false, $true, 1, $false, "542×9=", 2)
$d.Content.Find.Execute("554×9=", $true, $false, $false, $false, $false, $true, 1, $false, "809×8=", 2)
$d.Content.Find.Execute("192×5=", $true, $false, $false, $false, $false, $true, 1, $false, "905×2=", 2)
$d.Content.Find.Execute("412×5=", $true, $false, $false, $false, $false, $true, 1, $false, "268×5=", 2)
$d.Content.Find.Execute("283×4=", $true, $false, $false, $false, $false, $true, 1, $false, "761×8=", 2)
$d.Content.Find.Execute("265×5=", $true, $false, $false, $false, $false, $true, 1, $false, "565×7=", 2)
